# before.xlsx -> after.xlsx
#
# The sheet "1a" (a copy of a DNI list originally kept under
# files/identidades/secundaria) is being repurposed as the "dni_tutores"
# sheet under files/tutores. We:
#   1. Duplicate the sheet (Copy keeps styles/col widths/namespaces, and the
#      duplicate is assigned a fresh, incremented sheetId - matching the
#      sheetId="2" seen in the target workbook.xml).
#   2. Delete the original "1a" sheet, leaving only the duplicate.
#   3. Rename the duplicate to "dni_tutores".
#   4. Clear the number-format style that used to sit on A1 (41550112 is no
#      longer shown with a thousands separator).
#   5. Append the new tutor DNI row (A3 = 20771757).
#   6. Leave the selection where the editor ended up (B3), matching the
#      sheetView's final <selection> in the saved file.

$wb = $excel.ActiveWorkbook

# 1) Duplicate "1a" right after itself.
$source = $wb.Worksheets.Item("1a")
$source.Copy($null, $source)

# 2) Re-fetch the original sheet by name (object refs can go stale across
#    structural operations) and delete it, keeping only the duplicate.
$original = $wb.Worksheets.Item("1a")
$original.Delete()

# 3) The surviving duplicate is now named "1a (2)"; rename it.
$ws = $wb.Worksheets.Item("1a (2)")
$ws.Name = "dni_tutores"

# 4) A1 no longer carries the thousands-separator number format.
$ws.Range("A1").ClearFormats()

# 5) New tutor DNI.
$ws.Range("A3").Value = 20771757

# 6) Match the final selection left in the sheet.
$ws.Range("B3").Select()
